# Apply the "qry_Knights_ReleaseFish_EDI" update:
#  - rename the worksheet and the matching defined name
#  - extend the defined name's range to include the new last row
#  - append the new data row (2566) at the bottom of the table
#  - leave the active selection on C5, matching the authored edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (was "ReleaseFish___EDI_query")
$ws.Name = "qry_Knights_ReleaseFish_EDI"

# Append the new row of data at the end of the table
$ws.Cells.Item(2566, 1).Value = 9
$ws.Cells.Item(2566, 2).Value = 2565
$ws.Cells.Item(2566, 3).Value = 619

# Update the workbook-level defined name to track the sheet rename and
# the newly-extended data range
$definedName = $wb.Names.Item(1)
$definedName.Name = "qry_Knights_ReleaseFish_EDI"
$definedName.RefersTo = "=qry_Knights_ReleaseFish_EDI!`$A`$1:`$D`$2566"

# Move/restore the active selection to C5
$ws.Range("C5").Select()
